# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 0.581081081081081
$ws.Range("D3").Value = 0.9054054054054054
$ws.Range("H3").Value = 0.6533333333333333
$ws.Range("I3").Value = 0.09352971173170951
$ws.Range("J3").Value = 0.4864864864864865
$ws.Range("K3").Value = 121.7702702702703

$ws.Range("Q3").Value = 9
$ws.Range("S3").Value = 46
$ws.Range("T3").Value = 113
$ws.Range("U3").Value = 160
$ws.Range("V3").Value = 892
$ws.Range("X3").Value = 855
$ws.Range("Y3").Value = 788
$ws.Range("Z3").Value = 741

$ws.Range("AF3").Value = 0.990011
$ws.Range("AH3").Value = 0.948946
$ws.Range("AI3").Value = 0.874584
$ws.Range("AJ3").Value = 0.82242
